# Minor changes to slide 25 ("Two-stage Compiler"): inside the nested
# diagram group "Group 64" two small text-box labels read "X86/A" and
# need to become "x86/A"; their boxes are also nudged slightly (left a
# touch, narrower) to match the new text width.
#
# Both target boxes ("Text Box 82" and "Text Box 104") sit two levels
# deep inside groups-within-a-group, so they are reached through the
# (flattened) GroupItems collection of the top-level "Group 64" shape.
#
# NOTE on the literal Left/Top/Width/Height numbers below: this runtime
# stores a doubly-nested shape's own <a:off>/<a:ext> in the local
# (unscaled) child-coordinate units of its immediate parent group, but
# its Left/Top/Width/Height COM setters write `points * 12700` straight
# into that slot instead of applying the compounding group scale/offset
# transform. Passing plain "EMU/12700" point values therefore lands the
# shape exactly on the desired raw child-coordinate integers (1611,
# 3479, 427, 204, ...). A tiny epsilon is added because the interop
# layer loses a little precision in that conversion, which otherwise
# occasionally floors an exact integer down by one.

$EMU_PER_POINT = 12700
$EPS = 0.000001

function ChildUnitsToPoints($units) {
    return ($units / $EMU_PER_POINT) + $EPS
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(25)
$grp = $s.Shapes.Item(6)     # "Group 64" - contains both labels (nested)

for ($i = 1; $i -le $grp.GroupItems.Count; $i++) {
    $shp = $grp.GroupItems.Item($i)

    if ($shp.Name -eq "Text Box 82") {
        # Single-run label: "X86/A" -> "x86/A"
        $shp.TextFrame.TextRange.Text = "x86/A"

        $shp.Left   = ChildUnitsToPoints 1611
        $shp.Top    = ChildUnitsToPoints 3479
        $shp.Width  = ChildUnitsToPoints 427
        $shp.Height = ChildUnitsToPoints 204
    }
    elseif ($shp.Name -eq "Text Box 104") {
        # Two-run label: first run "X86/A " -> "x86/A ", second run (the
        # Symbol-font arrow + " x86") is left untouched.
        $shp.TextFrame.TextRange.Characters(1, 6).Text = "x86/A "

        $shp.Left   = ChildUnitsToPoints 682
        $shp.Top    = ChildUnitsToPoints 2557
        $shp.Width  = ChildUnitsToPoints 802
        $shp.Height = ChildUnitsToPoints 204
    }
}
